$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.897.93"
$ws.Range("E2").Value = "  -1.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.243.64"
$ws.Range("E3").Value = "  -1.66%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.41"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.09"
$ws.Range("E6").Value = "  -1.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.522"
$ws.Range("E7").Value = "  -1.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.487"
$ws.Range("E9").Value = "  -1.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.52"
$ws.Range("E10").Value = "  -4.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0809"
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.74"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.594.94"
$ws.Range("E14").Value = "  -1.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.33"
$ws.Range("E15").Value = "  -0.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.244.57"
$ws.Range("E16").Value = "  -2.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.778"
$ws.Range("E17").Value = "  -2.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.852.06"
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.09"
$ws.Range("E19").Value = "  -3.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0897"
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.90"
$ws.Range("E21").Value = "  -1.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.07"
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.80"
$ws.Range("E23").Value = "  -2.58%  "
$ws.Range("E24").Value = "  -1.63%  "
$ws.Range("E25").Value = "  -1.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("B27").Value = "InjectiveProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "37.70"
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.13"
$ws.Range("E28").Value = "  -3.25%  "
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.43"
$ws.Range("E30").Value = "  -1.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.30"
$ws.Range("E31").Value = "  +5.06%  "
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.13"
$ws.Range("E33").Value = "  -2.45%  "
$ws.Range("E34").Value = "  -2.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.35"
$ws.Range("E35").Value = "  +1.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0715"
$ws.Range("E36").Value = "  -3.54%  "
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("E39").Value = "  -3.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.79"
$ws.Range("E40").Value = "  -3.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.02"
$ws.Range("E41").Value = "  -2.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.935.25"
$ws.Range("E42").Value = "  -3.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0279"
$ws.Range("E43").Value = "  -2.03%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.50"
$ws.Range("E44").Value = "  -1.93%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.17"
$ws.Range("E45").Value = "  -10.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.87"
$ws.Range("E46").Value = "  -3.18%  "
$ws.Range("E47").Value = "  -3.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.59"
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.466.31"
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.88"
$ws.Range("E50").Value = "  -1.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "90.66"
$ws.Range("E51").Value = "  -1.66%  "
